$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2.jpg"
$ws.Range("A4").Value = "3.JPG"
$ws.Range("A6").Value = "5.JPG"
$ws.Range("A7").Value = "100.jpg"

$ws.Range("A3").Select()
